# 4.0.3 model and data
# Split the combined "BVTQaZ" and "VTQaZ" transportation CSV rows on the
# "Boolean" sheet into their per-mode-of-travel CSV files (LDVs, HDVs,
# aircraft, rail, ships, motorbikes), and update the active-sheet /
# selection state left behind by the editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Boolean" sheet: expand the two combined rows into six rows each.
# ---------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv". Make room for five
# more rows right below it, then fill A17:A22 with the split-out files.
$wsBool.Rows("18:22").Insert()
$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the first insert, "trans/VTQaZ/VTQaZ.csv" (formerly row 21) now
# sits at row 26. Make room for five more rows below it and fill
# A26:A31 with the split-out files.
$wsBool.Rows("27:31").Insert()
$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# "trans/VTStFES/VTStFES.csv" is now the last data row (row 32). Leave
# six blank formatted rows below it, matching the trailing whitespace
# left in the sheet by the editing session.
$wsBool.Rows("33:38").Insert()

# ---------------------------------------------------------------------
# View / selection state.
# ---------------------------------------------------------------------
# "Integer" sheet loses its tab selection but keeps a remembered
# selection at A13.
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Activate()
$wsInt.Range("A13").Select()

# "Boolean" sheet remembers a selection at its new last data row (A32).
$wsBool.Activate()
$wsBool.Range("A32").Select()

# "About" sheet becomes the active tab on reopen.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
